$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.801.67"
$ws.Range("E2").Value = "  -3.43%  "

$ws.Range("D3").Value = "2.912.11"
$ws.Range("E3").Value = "  -4.11%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.504"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.85%  "

$ws.Range("D9").Value = "2.912.03"
$ws.Range("E9").Value = "  -3.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.28%  "

$ws.Range("E11").Value = "  -4.85%  "

$ws.Range("E12").Value = "  -4.30%  "

$ws.Range("E13").Value = "  -4.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.17%  "

$ws.Range("E15").Value = "  +0.00%  "

$ws.Range("D16").Value = "3.392.00"
$ws.Range("E16").Value = "  -4.12%  "

$ws.Range("D17").Value = "60.735.64"
$ws.Range("E17").Value = "  -3.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.58%  "

$ws.Range("D19").Value = "2.908.74"
$ws.Range("E19").Value = "  -4.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.85%  "

$ws.Range("E22").Value = "  -2.26%  "

$ws.Range("E23").Value = "  -4.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.25%  "

$ws.Range("E26").Value = "  -5.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.17%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.64%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.45%  "

$ws.Range("E32").Value = "  -2.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.21%  "

$ws.Range("E34").Value = "  -4.31%  "

$ws.Range("D35").Value = "0.0₃0864"
$ws.Range("E35").Value = "  -0.63%  "

$ws.Range("E36").Value = "  -2.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.77%  "

$ws.Range("E39").Value = "  -0.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.65%  "

$ws.Range("E41").Value = "  -5.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.01%  "

$ws.Range("E43").Value = "  -4.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "376.16"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0348"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.36%  "

$ws.Range("D47").Value = "2.667.70"
$ws.Range("E47").Value = "  -2.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.33%  "

$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.75%  "

$ws.Range("E51").Value = "  -2.01%  "
